# Updated the backlog to reflect the product requirements.
# The priority has been updated to reflect what is needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the Priority column (A) for rows 3-15, matching the center-aligned
# numeric formatting already used by A2.
$templateCell = $ws.Range("A2")

$priorities = [ordered]@{
    "A3"  = 1
    "A4"  = 1
    "A5"  = 1
    "A6"  = 2
    "A7"  = 2
    "A8"  = 2
    "A9"  = 3
    "A10" = 3
    "A11" = 3
    "A12" = 4
    "A13" = 1
    "A14" = 4
    "A15" = 4
}

foreach ($addr in $priorities.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $priorities[$addr]
    $cell.HorizontalAlignment = $templateCell.HorizontalAlignment
}

# Update the active selection/view: select A15 only (previously the whole
# A1:D15 block was selected with B15 as the active cell, and B1 was pinned
# as the top-left visible cell).
$ws.Range("A15").Select()
